# Berechnung und Reports werden richtig erstellt
#
# The workbook lists 20 randomly-drawn students (Platznummer / Studentenname)
# in Tabelle1!B2:B21. The underlying calculation/report generation was fixed,
# which produced a different (corrected) random draw of students for the
# same 20 slots. Update the "Studentenname" column accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value  = "Student133"
$ws.Range("B3").Value  = "Student147"
$ws.Range("B4").Value  = "Student14"
$ws.Range("B5").Value  = "Student25"
$ws.Range("B6").Value  = "Student146"
$ws.Range("B7").Value  = "Student75"
$ws.Range("B8").Value  = "Student20"
$ws.Range("B9").Value  = "Student138"
$ws.Range("B10").Value = "Student69"
$ws.Range("B11").Value = "Student145"
$ws.Range("B12").Value = "Student122"
$ws.Range("B13").Value = "Student7"
$ws.Range("B14").Value = "Student40"
$ws.Range("B15").Value = "Student8"
$ws.Range("B16").Value = "Student140"
$ws.Range("B17").Value = "Student123"
$ws.Range("B18").Value = "Student143"
$ws.Range("B19").Value = "Student82"
$ws.Range("B20").Value = "Student120"
$ws.Range("B21").Value = "Student33"
